# Import section update: add dob/address/state/city/pincode/company/department/designation/others columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E: dob ---------------------------------------------------
$ws.Range("E1").Value = "dob"
$ws.Range("E2").Value = "23/1/2023"
$ws.Range("E3").Value = 44571
$ws.Range("E4").Value = "31/12/2022"
$ws.Range("E5").Value = 44927

# Apply a date number format to the whole dob column, reusing the same
# style record for every cell (copy/paste-format keeps the style table
# from growing one entry per cell).
$ws.Range("E2").NumberFormat = "mm-dd-yy"
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E3:E5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Column F: address ------------------------------------------------
$ws.Range("F1").Value = "address"
$ws.Range("F2").Value = "100/1 rajwada,raipur"
$ws.Range("F3").Value = "100/1 rajwada,raipur"
$ws.Range("F4").Value = "100/1 rajwada,raipur"
$ws.Range("F5").Value = "100/1 rajwada,raipur"

# --- Columns G:M headers ----------------------------------------------
$ws.Range("G1").Value = "state"
$ws.Range("H1").Value = "city"
$ws.Range("I1").Value = "pincode"
$ws.Range("J1").Value = "company"
$ws.Range("K1").Value = "department"
$ws.Range("L1").Value = "designation"
$ws.Range("M1").Value = "others"

# --- Columns G:M row 2 (first data values) -----------------------------
$ws.Range("G2").Value = "Chhattisgarh"
$ws.Range("H2").Value = "raipur"
$ws.Range("I2").Value = 12345
$ws.Range("J2").Value = "company name"
$ws.Range("K2").Value = "department name"
$ws.Range("L2").Value = "post"
$ws.Range("M2").Value = "abc"

# --- Columns G:M rows 3-5 (repeat the same values) ----------------------
foreach ($r in 3..5) {
    $ws.Range("G$r").Value = "Chhattisgarh"
    $ws.Range("H$r").Value = "raipur"
    $ws.Range("I$r").Value = 12345
    $ws.Range("J$r").Value = "company name"
    $ws.Range("K$r").Value = "department name"
    $ws.Range("L$r").Value = "post"
    $ws.Range("M$r").Value = "abc"
}

# --- Column widths (approximate the bestFit widths from the source file) --
$ws.Columns.Item(5).ColumnWidth = 10.7109375 - 0.8333333333333334
$ws.Columns.Item(6).ColumnWidth = 20 - 0.8333333333333334
$ws.Columns.Item(7).ColumnWidth = 12.140625 - 0.8333333333333334
$ws.Columns.Item(10).ColumnWidth = 14.5703125 - 0.8333333333333334
$ws.Columns.Item(11).ColumnWidth = 17.28515625 - 0.8333333333333334
$ws.Columns.Item(12).ColumnWidth = 11.42578125 - 0.8333333333333334

# --- Selection / view state --------------------------------------------
$ws.Range("M3:M5").Select() | Out-Null
